{"js": "// Actualizacion del manual de uso\n// Applies the textual changes described by the commit diff:\n//  1. \"cuatro\" -> \"cinco\" (number of screens)\n//  2. \"La segunda pantalla es el escenario de juego.\" ->\n//     \"La segunda pantalla es la introducci\u00f3n una peque\u00f1a introducci\u00f3n de\n//      la historia. La tercera es el escenario de juego.\"\n//  3. \"tercera\" (pause screen) -> \"cuarta\", and \"cuarta\" (final screen)\n//     -> \"quinta\" (renumbering the two last screens). These two replaces\n//     are done with search terms that sit entirely on one side of the\n//     \"_GoBack\" bookmark so the bookmark is not disturbed.\n//  4. Cosmetic run-merge around \"). La selecci\u00f3n de una \" (no visible text\n//     change, the two adjacent runs become one), matching the canonical\n//     OOXML produced by the edit.\n\nconst body = context.document.body;\n\n// 1) \"Existen cuatro pantallas.\" -> \"Existen cinco pantallas.\"\nlet results = body.search(\"Existen cuatro pantallas\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Existen cinco pantallas\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Expand the description of the second screen and introduce the third one.\nresults = body.search(\"La segunda pantalla es el escenario de juego.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"La segunda pantalla es la introducci\u00f3n una peque\u00f1a introducci\u00f3n de la historia. La tercera es el escenario de juego.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3a) The final screen becomes the \"quinta\" (do this before 3b, otherwise\n//     the new \"cuarta pantalla\" text created in 3b would also match here).\nresults = body.search(\"cuarta\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"quinta\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3b) The pause screen becomes the \"cuarta\" (this text lives entirely after\n//     the \"_GoBack\" bookmark, so the bookmark stays put).\nresults = body.search(\"tercera pantalla\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"cuarta pantalla\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) Normalize the adjacent runs around \"). La selecci\u00f3n de una \" (text is\n//    unchanged, but re-inserting merges the two runs into one, matching the\n//    canonical OOXML produced by the edit).\nresults = body.search(\"). La selecci\u00f3n de una \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"). La selecci\u00f3n de una \", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Actualizacion del manual de uso\n# Applies the textual changes described by the commit diff:\n#  1. \"cuatro\" -> \"cinco\" (number of screens)\n#  2. \"La segunda pantalla es el escenario de juego.\" ->\n#     \"La segunda pantalla es la introducci\u00f3n una peque\u00f1a introducci\u00f3n de\n#      la historia. La tercera es el escenario de juego.\"\n#  3. \"tercera\" (pause screen) -> \"cuarta\", and \"cuarta\" (final screen)\n#     -> \"quinta\" (renumbering the two last screens). These two replaces\n#     use search terms that sit entirely on one side of the \"_GoBack\"\n#     bookmark so the bookmark is not disturbed.\n#  4. Cosmetic run-merge around \"). La selecci\u00f3n de una \" (no visible text\n#     change, the two adjacent runs become one), matching the canonical\n#     OOXML produced by the edit.\n\n$d = $word.ActiveDocument\n\n# 1) \"Existen cuatro pantallas.\" -> \"Existen cinco pantallas.\"\n$rng = $d.Content\nif ($rng.Find.Execute(\"Existen cuatro pantallas\")) {\n    $rng.Text = \"Existen cinco pantallas\"\n}\n\n# 2) Expand the description of the second screen and introduce the third one.\n$rng = $d.Content\nif ($rng.Find.Execute(\"La segunda pantalla es el escenario de juego.\")) {\n    $rng.Text = \"La segunda pantalla es la introducci\u00f3n una peque\u00f1a introducci\u00f3n de la historia. La tercera es el escenario de juego.\"\n}\n\n# 3a) The final screen becomes the \"quinta\" (do this before 3b, otherwise the\n#     new \"cuarta pantalla\" text created in 3b would also match here).\n$rng = $d.Content\nif ($rng.Find.Execute(\"cuarta\")) {\n    $rng.Text = \"quinta\"\n}\n\n# 3b) The pause screen becomes the \"cuarta\" (this text lives entirely after\n#     the \"_GoBack\" bookmark, so the bookmark stays put).\n$rng = $d.Content\nif ($rng.Find.Execute(\"tercera pantalla\")) {\n    $rng.Text = \"cuarta pantalla\"\n}\n\n# 4) Normalize the adjacent runs around \"). La selecci\u00f3n de una \" (text is\n#    unchanged, but re-writing it merges the two runs into one, matching the\n#    canonical OOXML produced by the edit).\n$rng = $d.Content\nif ($rng.Find.Execute(\"). La selecci\u00f3n de una \")) {\n    $rng.Text = \"). La selecci\u00f3n de una \"\n}\n"}
